# Update the marksheet totals: B11 (Marking/Right count), B12 (Total), and
# E12 (the "Corr/total" display string) to reflect the new correct/total marks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 140
$ws.Range("E12").Value = "140/140"
